# Apply the "fixed string int error" update to the Correlation Table.
# - Column A: replace the coarse year-range labels with precise
#   year-month ranges.
# - Rows 12-16 (2020-01_2022-05) and 17-21 (2022-06_2023-12): the
#   Keyword/Correlation pairs were recomputed and reordered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: Year_Range relabeling (applies to all affected rows) ---
$yearRangeMap = @{
    "2008_2009" = "2007-12_2009-06"
    "2010_2019" = "2009-07_2019-12"
    "2020_2021" = "2020-01_2022-05"
    "2022_2023" = "2022-06_2023-12"
}

for ($r = 2; $r -le 21; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = $cell.Value2
    if ($yearRangeMap.ContainsKey($current)) {
        $cell.Value = $yearRangeMap[$current]
    }
}

# --- Rows 12-16: new Keyword order & recalculated Correlation values ---
$ws.Cells.Item(12, 2).Value = "inflation"
$ws.Cells.Item(12, 3).Value = -0.5714

$ws.Cells.Item(13, 2).Value = "uncertain"
$ws.Cells.Item(13, 3).Value = -0.4884

$ws.Cells.Item(14, 2).Value = "interest"
$ws.Cells.Item(14, 3).Value = -0.6074000000000001

$ws.Cells.Item(15, 2).Value = "invest"
$ws.Cells.Item(15, 3).Value = 0.1212

$ws.Cells.Item(16, 2).Value = "trade"
$ws.Cells.Item(16, 3).Value = -0.6751

# --- Rows 17-21: Keyword order unchanged, Correlation values recalculated ---
$ws.Cells.Item(17, 3).Value = -0.3442
$ws.Cells.Item(18, 3).Value = -0.2309
$ws.Cells.Item(19, 3).Value = -0.3358
$ws.Cells.Item(20, 3).Value = 0.6385
$ws.Cells.Item(21, 3).Value = -0.1606
